# QVA menu and its functionality added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: MSISDN becomes a plain numeric value (no text style), and the
# connection type narrows from "PREPAID (OR) POSTPAID" to "PREPAID".
$ws.Range("A2").Style = "Normal"
$ws.Range("A2").Value = 753011515
$ws.Range("F2").Value = "PREPAID"

# The trailing helper column G is no longer used on row 2.
$ws.Range("G2").Style = "Normal"
$ws.Range("G2").ClearContents()

# Row 3: MSISDN updated to a new number, stored as a true number (no style),
# and the connection type switches from PREPAID to POSTPAID (plain, no style).
$ws.Range("A3").Style = "Normal"
$ws.Range("A3").Value = 755843100
$ws.Range("F3").Style = "Normal"
$ws.Range("F3").Value = "POSTPAID"

# The trailing helper column G is no longer used on row 3 either.
$ws.Range("G3").Style = "Normal"
$ws.Range("G3").ClearContents()

# Row 4 (the second PREPAID/755841651 entry) is removed entirely.
$ws.Rows(4).Delete()

# Restore the saved selection to B2.
$ws.Range("B2").Select() | Out-Null
